# Change '_' to '-' in Excel metadata keyword names on the General_MD sheet.
# Close #5
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General_MD")

# Apply the edits in the same order the keywords appear top-to-bottom in the
# sheet (rows 25-28, then 32-35, then row 3), which is how a user manually
# correcting the sheet would naturally proceed.
$ws.Range("A25").Value = "SUBJECT-CODE"
$ws.Range("A26").Value = "SUBJECT-AREA_da"
$ws.Range("A27").Value = "SUBJECT-AREA_en"
$ws.Range("A28").Value = "SUBJECT-AREA_kl"
$ws.Range("A32").Value = "CREATION-DATE"
$ws.Range("A33").Value = "UPDATE-FREQUENCY"
$ws.Range("A34").Value = "LAST-UPDATED"
$ws.Range("A35").Value = "NEXT-UPDATE"
$ws.Range("A3").Value = "AXIS-VERSION"

# Leave the selection on A7, matching the saved view state.
$ws.Range("A7").Select()
